$wb = $excel.ActiveWorkbook

# --- Schedule sheet (sheet1): update cost + unit cost for the run ---
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Range("E2").Value = 726.52820175
$schedule.Range("F2").Value = 12.0127017485119

# --- Detailed sheet (sheet2): update per-interval Price values (and a couple of Type labels) ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B14").Value = 60.01196

$detailed.Range("B17").Value = 36.06
$detailed.Range("C17").Value = "historical"

$detailed.Range("B18").Value = 24.27246
$detailed.Range("C18").Value = "historical"

$detailed.Range("B19").Value = 0.7
$detailed.Range("B20").Value = -1.00223
$detailed.Range("B21").Value = -5.01
$detailed.Range("B22").Value = -5.58973
$detailed.Range("B23").Value = -5.01
$detailed.Range("B25").Value = -6.30983
$detailed.Range("B26").Value = -5.58973
$detailed.Range("B27").Value = -5.58973
$detailed.Range("B28").Value = -5.50985
$detailed.Range("B30").Value = -0.97989
$detailed.Range("B31").Value = 0.00003
$detailed.Range("B32").Value = 0.00948
$detailed.Range("B33").Value = 0.51
$detailed.Range("B34").Value = 0.51
$detailed.Range("B35").Value = -2.49263
$detailed.Range("B36").Value = -0.00598
$detailed.Range("B37").Value = 3.10484
$detailed.Range("B38").Value = 24.95242
$detailed.Range("B39").Value = 43.43343
$detailed.Range("B40").Value = 58.00487
$detailed.Range("B41").Value = 62.42689
$detailed.Range("B42").Value = 59.22102
$detailed.Range("B45").Value = 65
$detailed.Range("B46").Value = 62.19053
$detailed.Range("B47").Value = 62.33315
$detailed.Range("B48").Value = 60.36542
$detailed.Range("B49").Value = 63.5693
